$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap: row 109 gets old row 110 data
$ws.Cells.Item(109, 2).Value = 6830712
$ws.Cells.Item(109, 3).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(109, 4).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(109, 5).Value = 45201.91666666666
$ws.Cells.Item(109, 6).Value = "Tigres UANL Women"
$ws.Cells.Item(109, 7).Value = "Unam Pumas Women"
$ws.Cells.Item(109, 8).Value = 3
$ws.Cells.Item(109, 9).Value = 0
$ws.Cells.Item(109, 10).Value = "H"
$ws.Cells.Item(109, 11).Value = 1.181
$ws.Cells.Item(109, 12).Value = 6
$ws.Cells.Item(109, 13).Value = 10
$ws.Cells.Item(109, 14).Value = 1.1
$ws.Cells.Item(109, 15).Value = 9
$ws.Cells.Item(109, 16).Value = 19
$ws.Cells.Item(109, 17).Value = -2.5
$ws.Cells.Item(109, 18).Value = 1.8
$ws.Cells.Item(109, 19).Value = 2
$ws.Cells.Item(109, 20).Value = 3.75
$ws.Cells.Item(109, 21).Value = 1.8
$ws.Cells.Item(109, 22).Value = 2
$ws.Cells.Item(109, 23).Value = 0.1000000000000001
$ws.Cells.Item(109, 24).Value = -1
$ws.Cells.Item(109, 25).Value = -1
$ws.Cells.Item(109, 26).Value = 0.8
$ws.Cells.Item(109, 27).Value = -1
$ws.Cells.Item(109, 28).Value = -1
$ws.Cells.Item(109, 29).Value = 1

# Swap: row 110 gets old row 109 data
$ws.Cells.Item(110, 2).Value = 6830711
$ws.Cells.Item(110, 3).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(110, 4).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(110, 5).Value = 45201.91666666666
$ws.Cells.Item(110, 6).Value = "Atletico San Luis Women"
$ws.Cells.Item(110, 7).Value = "Tijuana Women"
$ws.Cells.Item(110, 8).Value = 0
$ws.Cells.Item(110, 9).Value = 2
$ws.Cells.Item(110, 10).Value = "A"
$ws.Cells.Item(110, 11).Value = 3
$ws.Cells.Item(110, 12).Value = 3.6
$ws.Cells.Item(110, 13).Value = 2
$ws.Cells.Item(110, 14).Value = 4
$ws.Cells.Item(110, 15).Value = 3.8
$ws.Cells.Item(110, 16).Value = 1.666
$ws.Cells.Item(110, 17).Value = 0.75
$ws.Cells.Item(110, 18).Value = 1.925
$ws.Cells.Item(110, 19).Value = 1.875
$ws.Cells.Item(110, 20).Value = 3
$ws.Cells.Item(110, 21).Value = 1.8
$ws.Cells.Item(110, 22).Value = 2
$ws.Cells.Item(110, 23).Value = -1
$ws.Cells.Item(110, 24).Value = -1
$ws.Cells.Item(110, 25).Value = 0.6659999999999999
$ws.Cells.Item(110, 26).Value = -1
$ws.Cells.Item(110, 27).Value = 0.875
$ws.Cells.Item(110, 28).Value = -1
$ws.Cells.Item(110, 29).Value = 1

# Swap: row 215 gets old row 216 data
$ws.Cells.Item(215, 2).Value = 7645772
$ws.Cells.Item(215, 3).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(215, 4).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(215, 5).Value = 45317.91666666666
$ws.Cells.Item(215, 6).Value = "Pachuca Women"
$ws.Cells.Item(215, 7).Value = "Queretaro Women"
$ws.Cells.Item(215, 8).Value = 4
$ws.Cells.Item(215, 9).Value = 1
$ws.Cells.Item(215, 10).Value = "H"
$ws.Cells.Item(215, 11).Value = 1.25
$ws.Cells.Item(215, 12).Value = 5.5
$ws.Cells.Item(215, 13).Value = 7.5
$ws.Cells.Item(215, 14).Value = 1.285
$ws.Cells.Item(215, 15).Value = 5.5
$ws.Cells.Item(215, 16).Value = 6.5
$ws.Cells.Item(215, 17).Value = -1.75
$ws.Cells.Item(215, 18).Value = 1.975
$ws.Cells.Item(215, 19).Value = 1.825
$ws.Cells.Item(215, 20).Value = 3.25
$ws.Cells.Item(215, 21).Value = 1.8
$ws.Cells.Item(215, 22).Value = 2
$ws.Cells.Item(215, 23).Value = 0.2849999999999999
$ws.Cells.Item(215, 24).Value = -1
$ws.Cells.Item(215, 25).Value = -1
$ws.Cells.Item(215, 26).Value = 0.9750000000000001
$ws.Cells.Item(215, 27).Value = -1
$ws.Cells.Item(215, 28).Value = 0.8
$ws.Cells.Item(215, 29).Value = -1

# Swap: row 216 gets old row 215 data
$ws.Cells.Item(216, 2).Value = 7645707
$ws.Cells.Item(216, 3).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(216, 4).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(216, 5).Value = 45317.91666666666
$ws.Cells.Item(216, 6).Value = "Tigres UANL Women"
$ws.Cells.Item(216, 7).Value = "Tijuana Women"
$ws.Cells.Item(216, 8).Value = 2
$ws.Cells.Item(216, 9).Value = 0
$ws.Cells.Item(216, 10).Value = "H"
$ws.Cells.Item(216, 11).Value = 1.125
$ws.Cells.Item(216, 12).Value = 7.5
$ws.Cells.Item(216, 13).Value = 15
$ws.Cells.Item(216, 14).Value = 1.166
$ws.Cells.Item(216, 15).Value = 7
$ws.Cells.Item(216, 16).Value = 15
$ws.Cells.Item(216, 17).Value = -2.25
$ws.Cells.Item(216, 18).Value = 1.825
$ws.Cells.Item(216, 19).Value = 1.975
$ws.Cells.Item(216, 20).Value = 3.75
$ws.Cells.Item(216, 21).Value = 1.85
$ws.Cells.Item(216, 22).Value = 1.95
$ws.Cells.Item(216, 23).Value = 0.1659999999999999
$ws.Cells.Item(216, 24).Value = -1
$ws.Cells.Item(216, 25).Value = -1
$ws.Cells.Item(216, 26).Value = -0.5
$ws.Cells.Item(216, 27).Value = 0.4875
$ws.Cells.Item(216, 28).Value = -1
$ws.Cells.Item(216, 29).Value = 0.95

# Swap: row 229 gets old row 231 data
$ws.Cells.Item(229, 2).Value = 7645781
$ws.Cells.Item(229, 3).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(229, 4).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(229, 5).Value = 45331.83333333334
$ws.Cells.Item(229, 6).Value = "Atletico San Luis Women"
$ws.Cells.Item(229, 7).Value = "Cruz Azul Women"
$ws.Cells.Item(229, 8).Value = 3
$ws.Cells.Item(229, 9).Value = 5
$ws.Cells.Item(229, 10).Value = "A"
$ws.Cells.Item(229, 11).Value = 2.1
$ws.Cells.Item(229, 12).Value = 3.6
$ws.Cells.Item(229, 13).Value = 2.8
$ws.Cells.Item(229, 14).Value = 2.45
$ws.Cells.Item(229, 15).Value = 3.6
$ws.Cells.Item(229, 16).Value = 2.375
$ws.Cells.Item(229, 17).Value = 0
$ws.Cells.Item(229, 18).Value = 1.95
$ws.Cells.Item(229, 19).Value = 1.85
$ws.Cells.Item(229, 20).Value = 2.75
$ws.Cells.Item(229, 21).Value = 1.75
$ws.Cells.Item(229, 22).Value = 1.95
$ws.Cells.Item(229, 23).Value = -1
$ws.Cells.Item(229, 24).Value = -1
$ws.Cells.Item(229, 25).Value = 1.375
$ws.Cells.Item(229, 26).Value = -1
$ws.Cells.Item(229, 27).Value = 0.8500000000000001
$ws.Cells.Item(229, 28).Value = 0.75
$ws.Cells.Item(229, 29).Value = -1

# Swap: row 231 gets old row 229 data
$ws.Cells.Item(231, 2).Value = 7645712
$ws.Cells.Item(231, 3).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(231, 4).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(231, 5).Value = 45331.83333333334
$ws.Cells.Item(231, 6).Value = "Chivas Guadalajara Women"
$ws.Cells.Item(231, 7).Value = "Santos Laguna Women"
$ws.Cells.Item(231, 8).Value = 10
$ws.Cells.Item(231, 9).Value = 2
$ws.Cells.Item(231, 10).Value = "H"
$ws.Cells.Item(231, 11).Value = 1.1
$ws.Cells.Item(231, 12).Value = 8
$ws.Cells.Item(231, 13).Value = 13
$ws.Cells.Item(231, 14).Value = 1.03
$ws.Cells.Item(231, 15).Value = 17
$ws.Cells.Item(231, 16).Value = 41
$ws.Cells.Item(231, 17).Value = -3.75
$ws.Cells.Item(231, 18).Value = 1.775
$ws.Cells.Item(231, 19).Value = 1.925
$ws.Cells.Item(231, 20).Value = 4.75
$ws.Cells.Item(231, 21).Value = 1.9
$ws.Cells.Item(231, 22).Value = 1.9
$ws.Cells.Item(231, 23).Value = 0.03000000000000003
$ws.Cells.Item(231, 24).Value = -1
$ws.Cells.Item(231, 25).Value = -1
$ws.Cells.Item(231, 26).Value = 0.7749999999999999
$ws.Cells.Item(231, 27).Value = -1
$ws.Cells.Item(231, 28).Value = 0.8999999999999999
$ws.Cells.Item(231, 29).Value = -1

# Swap: row 245 gets old row 246 data
$ws.Cells.Item(245, 2).Value = 7645793
$ws.Cells.Item(245, 3).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(245, 4).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(245, 5).Value = 45355.83333333334
$ws.Cells.Item(245, 6).Value = "Queretaro Women"
$ws.Cells.Item(245, 7).Value = "Cruz Azul Women"
$ws.Cells.Item(245, 8).Value = 3
$ws.Cells.Item(245, 9).Value = 0
$ws.Cells.Item(245, 10).Value = "H"
$ws.Cells.Item(245, 11).Value = 1.909
$ws.Cells.Item(245, 12).Value = 3.5
$ws.Cells.Item(245, 13).Value = 3.3
$ws.Cells.Item(245, 14).Value = 2.1
$ws.Cells.Item(245, 15).Value = 3.4
$ws.Cells.Item(245, 16).Value = 2.875
$ws.Cells.Item(245, 17).Value = -0.25
$ws.Cells.Item(245, 18).Value = 1.875
$ws.Cells.Item(245, 19).Value = 1.925
$ws.Cells.Item(245, 20).Value = 2.75
$ws.Cells.Item(245, 21).Value = 2
$ws.Cells.Item(245, 22).Value = 1.8
$ws.Cells.Item(245, 23).Value = 1.1
$ws.Cells.Item(245, 24).Value = -1
$ws.Cells.Item(245, 25).Value = -1
$ws.Cells.Item(245, 26).Value = 0.875
$ws.Cells.Item(245, 27).Value = -1
$ws.Cells.Item(245, 28).Value = 0.5
$ws.Cells.Item(245, 29).Value = -0.5

# Swap: row 246 gets old row 245 data
$ws.Cells.Item(246, 2).Value = 7645794
$ws.Cells.Item(246, 3).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(246, 4).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(246, 5).Value = 45355.83333333334
$ws.Cells.Item(246, 6).Value = "Toluca Women"
$ws.Cells.Item(246, 7).Value = "Club America Women"
$ws.Cells.Item(246, 8).Value = 3
$ws.Cells.Item(246, 9).Value = 0
$ws.Cells.Item(246, 10).Value = "H"
$ws.Cells.Item(246, 11).Value = 6
$ws.Cells.Item(246, 12).Value = 5
$ws.Cells.Item(246, 13).Value = 1.333
$ws.Cells.Item(246, 14).Value = 5.75
$ws.Cells.Item(246, 15).Value = 5
$ws.Cells.Item(246, 16).Value = 1.363
$ws.Cells.Item(246, 17).Value = 1.5
$ws.Cells.Item(246, 18).Value = 1.825
$ws.Cells.Item(246, 19).Value = 1.975
$ws.Cells.Item(246, 20).Value = 3.25
$ws.Cells.Item(246, 21).Value = 1.825
$ws.Cells.Item(246, 22).Value = 1.975
$ws.Cells.Item(246, 23).Value = 4.75
$ws.Cells.Item(246, 24).Value = -1
$ws.Cells.Item(246, 25).Value = -1
$ws.Cells.Item(246, 26).Value = 0.825
$ws.Cells.Item(246, 27).Value = -1
$ws.Cells.Item(246, 28).Value = -0.5
$ws.Cells.Item(246, 29).Value = 0.4875

# Swap: row 251 gets old row 252 data
$ws.Cells.Item(251, 2).Value = 7645798
$ws.Cells.Item(251, 3).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(251, 4).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(251, 5).Value = 45359.83333333334
$ws.Cells.Item(251, 6).Value = "Atlas Women"
$ws.Cells.Item(251, 7).Value = "Unam Pumas Women"
$ws.Cells.Item(251, 8).Value = 0
$ws.Cells.Item(251, 9).Value = 2
$ws.Cells.Item(251, 10).Value = "A"
$ws.Cells.Item(251, 11).Value = 2.4
$ws.Cells.Item(251, 12).Value = 3.6
$ws.Cells.Item(251, 13).Value = 2.4
$ws.Cells.Item(251, 14).Value = 2.375
$ws.Cells.Item(251, 15).Value = 3.75
$ws.Cells.Item(251, 16).Value = 2.375
$ws.Cells.Item(251, 17).Value = 0
$ws.Cells.Item(251, 18).Value = 1.925
$ws.Cells.Item(251, 19).Value = 1.875
$ws.Cells.Item(251, 20).Value = 3
$ws.Cells.Item(251, 21).Value = 1.825
$ws.Cells.Item(251, 22).Value = 1.975
$ws.Cells.Item(251, 23).Value = -1
$ws.Cells.Item(251, 24).Value = -1
$ws.Cells.Item(251, 25).Value = 1.375
$ws.Cells.Item(251, 26).Value = -1
$ws.Cells.Item(251, 27).Value = 0.875
$ws.Cells.Item(251, 28).Value = -1
$ws.Cells.Item(251, 29).Value = 0.9750000000000001

# Swap: row 252 gets old row 251 data
$ws.Cells.Item(252, 2).Value = 7926076
$ws.Cells.Item(252, 3).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(252, 4).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(252, 5).Value = 45359.83333333334
$ws.Cells.Item(252, 6).Value = "Leon Women"
$ws.Cells.Item(252, 7).Value = "Monterrey Women"
$ws.Cells.Item(252, 8).Value = 1
$ws.Cells.Item(252, 9).Value = 4
$ws.Cells.Item(252, 10).Value = "A"
$ws.Cells.Item(252, 11).Value = 4.5
$ws.Cells.Item(252, 12).Value = 4
$ws.Cells.Item(252, 13).Value = 1.571
$ws.Cells.Item(252, 14).Value = 4.75
$ws.Cells.Item(252, 15).Value = 4
$ws.Cells.Item(252, 16).Value = 1.533
$ws.Cells.Item(252, 17).Value = 1
$ws.Cells.Item(252, 18).Value = 1.825
$ws.Cells.Item(252, 19).Value = 1.975
$ws.Cells.Item(252, 20).Value = 3
$ws.Cells.Item(252, 21).Value = 1.975
$ws.Cells.Item(252, 22).Value = 1.825
$ws.Cells.Item(252, 23).Value = -1
$ws.Cells.Item(252, 24).Value = -1
$ws.Cells.Item(252, 25).Value = 0.5329999999999999
$ws.Cells.Item(252, 26).Value = -1
$ws.Cells.Item(252, 27).Value = 0.9750000000000001
$ws.Cells.Item(252, 28).Value = 0.9750000000000001
$ws.Cells.Item(252, 29).Value = -1

# Rotate: row 263 gets old row 265 data
$ws.Cells.Item(263, 2).Value = 7645804
$ws.Cells.Item(263, 3).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(263, 4).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(263, 5).Value = 45369.83333333334
$ws.Cells.Item(263, 6).Value = "Puebla Women"
$ws.Cells.Item(263, 7).Value = "Mazatlan FC Women"
$ws.Cells.Item(263, 8).Value = 1
$ws.Cells.Item(263, 9).Value = 2
$ws.Cells.Item(263, 10).Value = "A"
$ws.Cells.Item(263, 11).Value = 1.666
$ws.Cells.Item(263, 12).Value = 4
$ws.Cells.Item(263, 13).Value = 3.8
$ws.Cells.Item(263, 14).Value = 1.333
$ws.Cells.Item(263, 15).Value = 4.5
$ws.Cells.Item(263, 16).Value = 7.5
$ws.Cells.Item(263, 17).Value = -1.5
$ws.Cells.Item(263, 18).Value = 1.925
$ws.Cells.Item(263, 19).Value = 1.875
$ws.Cells.Item(263, 20).Value = 3
$ws.Cells.Item(263, 21).Value = 1.75
$ws.Cells.Item(263, 22).Value = 2.05
$ws.Cells.Item(263, 23).Value = -1
$ws.Cells.Item(263, 24).Value = -1
$ws.Cells.Item(263, 25).Value = 6.5
$ws.Cells.Item(263, 26).Value = -1
$ws.Cells.Item(263, 27).Value = 0.875
$ws.Cells.Item(263, 28).Value = 0
$ws.Cells.Item(263, 29).Value = -0

# Rotate: row 264 gets old row 263 data
$ws.Cells.Item(264, 2).Value = 7645807
$ws.Cells.Item(264, 3).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(264, 4).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(264, 5).Value = 45369.83333333334
$ws.Cells.Item(264, 6).Value = "Club Necaxa Women"
$ws.Cells.Item(264, 7).Value = "Leon Women"
$ws.Cells.Item(264, 8).Value = 2
$ws.Cells.Item(264, 9).Value = 1
$ws.Cells.Item(264, 10).Value = "H"
$ws.Cells.Item(264, 11).Value = 4.333
$ws.Cells.Item(264, 12).Value = 4
$ws.Cells.Item(264, 13).Value = 1.571
$ws.Cells.Item(264, 14).Value = 7
$ws.Cells.Item(264, 15).Value = 4.2
$ws.Cells.Item(264, 16).Value = 1.363
$ws.Cells.Item(264, 17).Value = 1.5
$ws.Cells.Item(264, 18).Value = 1.75
$ws.Cells.Item(264, 19).Value = 1.95
$ws.Cells.Item(264, 20).Value = 2.75
$ws.Cells.Item(264, 21).Value = 1.8
$ws.Cells.Item(264, 22).Value = 2
$ws.Cells.Item(264, 23).Value = 6
$ws.Cells.Item(264, 24).Value = -1
$ws.Cells.Item(264, 25).Value = -1
$ws.Cells.Item(264, 26).Value = 0.75
$ws.Cells.Item(264, 27).Value = -1
$ws.Cells.Item(264, 28).Value = 0.4
$ws.Cells.Item(264, 29).Value = -0.5

# Rotate: row 265 gets old row 264 data
$ws.Cells.Item(265, 2).Value = 7645806
$ws.Cells.Item(265, 3).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(265, 4).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(265, 5).Value = 45369.83333333334
$ws.Cells.Item(265, 6).Value = "Atletico San Luis Women"
$ws.Cells.Item(265, 7).Value = "Atlas Women"
$ws.Cells.Item(265, 8).Value = 0
$ws.Cells.Item(265, 9).Value = 0
$ws.Cells.Item(265, 10).Value = "D"
$ws.Cells.Item(265, 11).Value = 3.2
$ws.Cells.Item(265, 12).Value = 3.6
$ws.Cells.Item(265, 13).Value = 1.909
$ws.Cells.Item(265, 14).Value = 3.1
$ws.Cells.Item(265, 15).Value = 3.6
$ws.Cells.Item(265, 16).Value = 2
$ws.Cells.Item(265, 17).Value = 0.25
$ws.Cells.Item(265, 18).Value = 1.975
$ws.Cells.Item(265, 19).Value = 1.825
$ws.Cells.Item(265, 20).Value = 3
$ws.Cells.Item(265, 21).Value = 1.9
$ws.Cells.Item(265, 22).Value = 1.9
$ws.Cells.Item(265, 23).Value = -1
$ws.Cells.Item(265, 24).Value = 2.6
$ws.Cells.Item(265, 25).Value = -1
$ws.Cells.Item(265, 26).Value = 0.4875
$ws.Cells.Item(265, 27).Value = -0.5
$ws.Cells.Item(265, 28).Value = -1
$ws.Cells.Item(265, 29).Value = 0.8999999999999999

# Direct odds correction: row 278 (R:S:T:U:V)
$ws.Cells.Item(278, 18).Value = 1.825
$ws.Cells.Item(278, 19).Value = 1.975
$ws.Cells.Item(278, 20).Value = 4
$ws.Cells.Item(278, 21).Value = 2
$ws.Cells.Item(278, 22).Value = 1.8

# Direct odds correction: row 281 (N:O:P:R:S)
$ws.Cells.Item(281, 14).Value = 6.5
$ws.Cells.Item(281, 15).Value = 4
$ws.Cells.Item(281, 16).Value = 1.4
$ws.Cells.Item(281, 18).Value = 1.975
$ws.Cells.Item(281, 19).Value = 1.825
